$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without letting Excel's
# auto-detection reinterpret a date-shaped string ("01/01/2023") as a
# serial date number (which would also mint a brand-new number-format
# style). We stage the text (quote-prefixed, so it is stored as text) in
# a scratch cell that already uses the destination's style, copy it, and
# paste-special "values only" into the destination -- this carries the
# text over while leaving the destination's existing style (s=2 / s=3)
# completely untouched. The scratch cell's original content is restored
# afterwards.
function Set-PlainText($addr, $text) {
    $scratch = $ws.Range("B9")
    $origFormula = $scratch.Formula
    $scratch.Value = "'" + $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Formula = $origFormula
}

# Ativação / Programa date bump: 01/01/2020 -> 01/01/2023
Set-PlainText "B8" "01/01/2023"
Set-PlainText "C8" "01/01/2023"
Set-PlainText "B15" "01/01/2023"
Set-PlainText "C15" "01/01/2023"

# Objectives (English) -- row 11
$ws.Range("B11").Value = "To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others"
$ws.Range("C11").Value = "To present concepts about renewable sources for the generation of thermal, electrical and vehicular energy, among others"

# Short syllabus (English) -- row 14
$ws.Range("B14").Value = "Renewable sources and clean technologies for energy generation. Study of current national and global systems."
$ws.Range("C14").Value = "Renewable sources and clean technologies for energy generation. Study of current national and global systems."

# Syllabus (English) -- row 16
$ws.Range("B16").Value = "National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles"
$ws.Range("C16").Value = "National and global energy systems: renewable and fossil sources. Energy generation from renewable sources: solar thermal and photovoltaic; wind; maritime. Generation of biomass for energy purposes. Management of urban solid waste: recyclable and non-recyclable; enterprise programs for reverse logistics; the issue of polymers; reforestation; processing of domestic wet waste. Integration of renewable sources for energy generation: hybrid thermal cycles"
